# Add the "undo/redo" design-decision bullet right after the paragraph
# that ends with "...using a dynamic approach that doesn't need to be
# changed upon the addition of new plants." (the last bullet of the
# "Design Changes (M2->M3)" list, numId=5), and before the "What Smells:"
# heading paragraph.

$d = $word.ActiveDocument

# 1. Locate the anchor paragraph (last numId=5 bullet before "What Smells:")
$anchorText = "How the Plants are dealt with in the controllers and how plant buttons were created has also been changed during this milestone, using a dynamic approach that doesn" + [char]0x2019 + "t need to be changed upon the addition of new plants. "

$found = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("How the Plants are dealt with")) {
        $found = $p
        break
    }
}

if ($found -eq $null) {
    throw "Anchor paragraph not found"
}

# 2. Insert a new paragraph right after it, inheriting the list formatting.
$newRange = $found.Range.InsertParagraphAfter()

# 3. The newly created paragraph is the one that now sits between the
#    anchor and what used to be the next paragraph. Grab it via Next.
$newPara = $found.Next()

# Apply the same list style / numbering as the anchor paragraph.
$newPara.Range.ListFormat.ApplyListTemplateWithLevel($found.Range.ListFormat.ListTemplate, $false, 2, $false, $false)
$newPara.Style = "ListParagraph"

# 4. Put the text into the new paragraph's range.
$undoText = "The undo redo feature was added using a doubly linked list, storing past plant placements, as well as future plant placements(Which only exist after undoing). Whenever a plant was placed, the future plants were cleared, and the plant was added to the past placements. At the end of each turn, the entire linked list was cleared."
$newPara.Range.Text = $undoText

# 5. Word keeps a single, auto-managed "_GoBack" bookmark that tracks the
#    location of the most recent edit; remove any pre-existing one and add
#    it at the end of the freshly-typed text (mirrors what Word does after
#    you type new content).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $newPara.Range)
